$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the three changed data cells in row 3
$ws.Range("E3").Value = 5
$ws.Range("G3").Value = -3
$ws.Range("H3").Value = 13

# Reflect the new active cell / selection saved with the sheet
$ws.Range("E3").Select()
